# RingampSurvey.xlsx edit:
#  - add a new 2024 paper row (row 135) to the "Ringamp Publication List" sheet
#  - annotate "Sim only" in the Notes column (D) for a number of existing rows
#  - the Analysis sheet's COUNTIF/SUM formulas and the charts that are driven
#    by it recompute automatically; chart3's plotted category range shrinks
#    by one row (the "All other" bucket is dropped from the bar chart)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ringamp Publication List")
$ws2 = $wb.Worksheets.Item("Analysis")

# ---------------------------------------------------------------------------
# 1. New row 135 - "A 12b 400MS/s 4-Time Interleaved Pipelined-SAR ADC..."
#    Field order below matches the order the strings were first authored in
#    (title, authors, doi, abstract, venue, notes) so shared-string indices
#    line up the same way they would from typing the row into Excel.
# ---------------------------------------------------------------------------
$ws1.Range("E135").Value = 'A 12b 400MS/s 4-Time Interleaved Pipelined-SAR ADC with Fully Differential Bias-enhanced Ring Amplifier'
$ws1.Range("G135").Value = 'Jian, Mingchao and Kong, Xiangjian and Zheng, Jiwei and Xie, Huanlin and Guo, Chunbing'
$ws1.Range("F135").Value = '10.1109/ICICM63644.2024.10814121'
$ws1.Range("H135").Value = 'A 12-bit 400MS/s 4-channel time-interleaved pipelined-SAR ADC with a fully differential bias-enhanced ring amplifier for low-power and high-speed application is presented. The sub-ADC adopts a two-stage pipeline quantization scheme with a 6-bit SAR ADC in the first stage and a 7-bit SAR ADC in the second stage with 1-hit overlapping. To further optimize the speed of the sub-ADCs, a bias-enhanced ring amplifier is proposed. A prototype ADC is designed and simulated in 65nm CMOS technology with a standard 1.2 V supply voltage. With the least-mean-square algorithm calibration, this ADC achieves SNDR of 64.5 dB and SFDR of 86.6 dB. The whole ADC consumes 10.3mW, achieving a FoM of 19.9 fJ/conv-step.'
$ws1.Range("B135").Value = 'ICICM'
$ws1.Range("D135").Value = 'Interleaved, Sim only'
$ws1.Range("A135").Value = 2024
$ws1.Range("C135").Value = 'Pipelined-SAR ADC'

# ---------------------------------------------------------------------------
# 2. Annotate "Sim only" in column D (Notes) for the rows that need it.
# ---------------------------------------------------------------------------
$simOnlyRows = @(7, 11, 13, 19, 20, 24, 26, 27, 28, 29, 36, 37, 61, 63, 85, 95, 101, 116, 132, 133)
foreach ($r in $simOnlyRows) {
    $ws1.Range("D$r").Value = 'Sim only'
}

# Row 124 already has a Notes value of "FIA" -> becomes "FIA, Sim only"
$ws1.Range("D124").Value = 'FIA, Sim only'

# ---------------------------------------------------------------------------
# 3. Chart "Publication Venue of Ringamp Papers" (3rd chart on Analysis):
#    its plotted range shrinks from Analysis!$A$38:$A$49 to $A$38:$A$48,
#    i.e. it no longer plots the trailing "All other" bucket.
# ---------------------------------------------------------------------------
$chartObjs = $ws2.ChartObjects()
$venueChart = $chartObjs.Item(3).Chart
$venueSeries = $venueChart.SeriesCollection(1)
$venueSeries.Formula = "=SERIES(,Analysis!`$A`$38:`$A`$48,Analysis!`$B`$38:`$B`$48,1)"

# ---------------------------------------------------------------------------
# 4. View-state bookkeeping (best effort): keep the publication list the
#    active sheet, scrolled down near the newly-added row, and leave the
#    Analysis sheet's selection on B24.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("B24").Select()

$ws1.Activate()
$ws1.Range("A74").Select()
